$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New cell values (plain, no special style) ---
$ws.Range("I4").Value = 6
$ws.Range("J4").Value = 7
$ws.Range("K4").Value = 8
$ws.Range("L4").Value = 9
$ws.Range("M4").Value = 10
$ws.Range("N4").Value = 11
$ws.Range("O4").Value = 12
$ws.Range("P4").Value = 13
$ws.Range("Q4").Value = 14
$ws.Range("R4").Value = 15
$ws.Range("S4").Value = 16
$ws.Range("T4").Value = 17

$ws.Range("K6").Value = 8

$ws.Range("K7").Value = 8
$ws.Range("L7").Value = 9
$ws.Range("M7").Value = 10
$ws.Range("N7").Value = 11
$ws.Range("O7").Value = 12
$ws.Range("P7").Value = 13

$ws.Range("U9").Value = 18
$ws.Range("V9").Value = 19
$ws.Range("W9").Value = 20
$ws.Range("X9").Value = 21
$ws.Range("Y9").Value = 22
$ws.Range("Z9").Value = 23

$ws.Range("L10").Value = 9
$ws.Range("M10").Value = 10
$ws.Range("N10").Value = 11
$ws.Range("O10").Value = 12

$ws.Range("AA12").Value = 24
$ws.Range("AB12").Value = 25
$ws.Range("AC12").Value = 26
$ws.Range("AD12").Value = 27
$ws.Range("AE12").Value = 28
$ws.Range("AF12").Value = 29

$ws.Range("P13").Value = 13
$ws.Range("Q13").Value = 14
$ws.Range("R13").Value = 15
$ws.Range("S13").Value = 16
$ws.Range("T13").Value = 17
$ws.Range("U13").Value = 18

# --- New cell values with centered alignment (no border) ---
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 2
$ws.Range("E2").HorizontalAlignment = -4108

$ws.Range("F3").Value = 3
$ws.Range("G3").Value = 4
$ws.Range("H3").Value = 5
$ws.Range("F3:H3").HorizontalAlignment = -4108

# --- New cell values with centered alignment + left/right thin border (per-cell) ---
$rowsWithBorderedValues = @(
  @{ Row = 5;  Cells = @{ D=1; E=2; F=3; G=4; H=5; I=6; J=7 } },
  @{ Row = 8;  Cells = @{ D=1; E=2; F=3; G=4; H=5; I=6 } },
  @{ Row = 11; Cells = @{ D=1; E=2; F=3; G=4; H=5; I=6; J=7; K=8; L=9; M=10; N=11 } }
)

foreach ($grp in $rowsWithBorderedValues) {
  $r = $grp.Row
  foreach ($col in $grp.Cells.Keys) {
    $addr = "$col$r"
    $cell = $ws.Range($addr)
    $cell.Value = $grp.Cells[$col]
    $cell.HorizontalAlignment = -4108
    $cell.Borders.Item(7).LineStyle = 1
    $cell.Borders.Item(10).LineStyle = 1
  }
}

# --- Sheet view selection ---
$ws.Range("T29").Select()

# --- Column widths (closest achievable approximation) ---
$ws.Columns.Item(1).ColumnWidth = 15
$ws.Columns.Item(2).ColumnWidth = 15.6666666667
$ws.Columns.Item(3).ColumnWidth = 14.3333333333

# --- Page setup (portrait) ---
$ws.PageSetup.Orientation = 1

Write-Host "done"
